# Update cryptocurrency Price (D) and Volume(1h) (E) columns to reflect refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''56.425.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.05%  '
$ws.Range("D3").Value = '''2.970.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.25%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''495.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.36%  '
$ws.Range("D6").Value = '''134.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''2.967.83'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("D9").Value = '''0.425'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.94%  '
$ws.Range("D10").Value = '''7.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("E11").Value = '  -3.78%  '
$ws.Range("E12").Value = '  -7.38%  '
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").Value = '''3.481.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.17%  '
$ws.Range("D15").Value = '''25.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("D16").Value = '''56.410.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.05%  '
$ws.Range("D17").Value = '''2.974.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.07%  '
$ws.Range("D18").Value = '''0.0000145'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.57%  '
$ws.Range("D19").Value = '''5.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").Value = '''12.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.94%  '
$ws.Range("D21").Value = '''7.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.73%  '
$ws.Range("D22").Value = '''324.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.46%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -7.86%  '
$ws.Range("D25").Value = '''61.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.51%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = '''0.162'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.85%  '
$ws.Range("E28").Value = '  -5.87%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '''6.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.01%  '
$ws.Range("D31").Value = '''6.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.72%  '
$ws.Range("D32").Value = '''1.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.72%  '
$ws.Range("E33").Value = '  -6.51%  '
$ws.Range("D34").Value = '''20.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.59%  '
$ws.Range("D35").Value = '''152.89'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").Value = '''4.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.93%  '
$ws.Range("D37").Value = '''1.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.66%  '
$ws.Range("D38").Value = '''5.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.43%  '
$ws.Range("D39").Value = '''0.0668'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.30%  '
$ws.Range("D40").Value = '''23.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").Value = '''3.005.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.08%  '
$ws.Range("E42").Value = '  -9.60%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '''0.637'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.76%  '
$ws.Range("D45").Value = '''0.989'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.27%  '
$ws.Range("D46").Value = '''1.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("D47").Value = '''2.204.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("D48").Value = '''3.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.82%  '
$ws.Range("E49").Value = '  +4.25%  '
$ws.Range("D50").Value = '''0.0237'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").Value = '''5.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.16%  '
